$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new data row (row 37) that was previously blank
$ws.Range("A37").Value = 45992
$ws.Range("B37").Value = 663
$ws.Range("C37").Value = 23
$ws.Range("D37").Value = 640

# Update the active selection to reflect the newly entered row
$ws.Range("A37:D37").Select()
